$d = $word.ActiveDocument

$d.Content.Find.Execute("68÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷2=", 2) | Out-Null
$d.Content.Find.Execute("95÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=", 2) | Out-Null
$d.Content.Find.Execute("96÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷7=", 2) | Out-Null
$d.Content.Find.Execute("77÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=", 2) | Out-Null
$d.Content.Find.Execute("53÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("60÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷5=", 2) | Out-Null
$d.Content.Find.Execute("27÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=", 2) | Out-Null
$d.Content.Find.Execute("30÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=", 2) | Out-Null
$d.Content.Find.Execute("24÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=", 2) | Out-Null
$d.Content.Find.Execute("65÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=", 2) | Out-Null
$d.Content.Find.Execute("51÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷9=", 2) | Out-Null
$d.Content.Find.Execute("41÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷2=", 2) | Out-Null
$d.Content.Find.Execute("91÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷4=", 2) | Out-Null
$d.Content.Find.Execute("34÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷8=", 2) | Out-Null
$d.Content.Find.Execute("76÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷4=", 2) | Out-Null
$d.Content.Find.Execute("90÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2) | Out-Null
$d.Content.Find.Execute("97÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=", 2) | Out-Null
$d.Content.Find.Execute("53÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=", 2) | Out-Null
$d.Content.Find.Execute("74÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷7=", 2) | Out-Null
$d.Content.Find.Execute("88÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷5=", 2) | Out-Null
$d.Content.Find.Execute("99÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷2=", 2) | Out-Null
$d.Content.Find.Execute("37÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷9=", 2) | Out-Null
$d.Content.Find.Execute("37÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=", 2) | Out-Null
$d.Content.Find.Execute("38÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=", 2) | Out-Null
$d.Content.Find.Execute("87÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷3=", 2) | Out-Null
